$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added at the top of the data block
# (row 805), pushing every existing record down by one row. Insert a
# fresh row at 805 (this shifts rows 805..833 down to 806..834, and the
# sheet's used range grows from A1:R833 to A1:R834 automatically).
$ws.Rows.Item(805).Insert()

# Populate the newly inserted row with the new observation. Everything
# except the date / volume / price columns repeats the values of the
# record that is now at row 806 (Camote, "1a (guarda)", Región del Maule).
$ws.Cells.Item(805, 1).Value = 10
$ws.Cells.Item(805, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(805, 3).Value = "La Araucanía"
$ws.Cells.Item(805, 4).Value = 45075
$ws.Cells.Item(805, 5).Value = 9
$ws.Cells.Item(805, 6).Value = 100112045
$ws.Cells.Item(805, 7).Value = "Zapallo"
$ws.Cells.Item(805, 8).Value = "Camote"
$ws.Cells.Item(805, 9).Value = "1a (guarda)"
$ws.Cells.Item(805, 10).Value = 680
$ws.Cells.Item(805, 11).Value = 500
$ws.Cells.Item(805, 12).Value = 500
$ws.Cells.Item(805, 13).Value = 500
$ws.Cells.Item(805, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(805, 15).Value = "Región del Maule"
$ws.Cells.Item(805, 16).Value = 500
$ws.Cells.Item(805, 17).Value = 1
$ws.Cells.Item(805, 18).Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of the
# column (it normally survives the row insert, but set it explicitly to
# be safe).
$ws.Cells.Item(805, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
